# Automatic update of files.
#
# The rows 16-38 (data rows) in the "Artfynd" sheet need to have their entire
# content permuted: each row's full set of column values (A..AY) is replaced
# by the values that originally belonged to a different row in that same
# block. Row 28 is left untouched.
#
# Mapping: new row -> source row (row that currently holds the content that
# should end up in "new row" after the edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$firstRow = 16
$lastRow = 38
$firstCol = 1   # A
$lastCol = 51   # AY

# new row number -> original row number supplying its content
$map = @{
    16 = 33
    17 = 16
    18 = 23
    19 = 18
    20 = 17
    21 = 29
    22 = 20
    23 = 30
    24 = 32
    25 = 19
    26 = 25
    27 = 35
    29 = 26
    30 = 34
    31 = 27
    32 = 22
    33 = 36
    34 = 37
    35 = 21
    36 = 24
    37 = 38
    38 = 31
}

# 1) Read the full original block (values only) into memory so that writing
#    the permuted data back doesn't clobber a row before it has been read.
#    Excel returns this Value2 array as a 1-based (COM SAFEARRAY) array.
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$srcValues = $srcRange.Value2

$rowCount = $lastRow - $firstRow + 1
$colCount = $lastCol - $firstCol + 1

# 2) Build the destination array applying the permutation; rows that are not
#    remapped (e.g. row 28) keep their original content.
#    New-Object arrays are 0-based, so keep that indexing explicit below.
$destValues = New-Object 'object[,]' $rowCount, $colCount

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $destIdx0 = $r - $firstRow
    if ($map.ContainsKey($r)) {
        $srcRow = $map[$r]
    } else {
        $srcRow = $r
    }
    $srcIdx1 = $srcRow - $firstRow + 1
    for ($c = 1; $c -le $colCount; $c++) {
        $destValues[$destIdx0, $c - 1] = $srcValues[$srcIdx1, $c]
    }
}

# 3) Write the permuted block back, skipping rows that are not remapped
#    (their content must stay completely untouched). String values are
#    written with a leading apostrophe to stop Excel from "smart"
#    re-interpreting them (e.g. turning a plain text date like
#    "2023-08-24" into a real date value) - the apostrophe itself is not
#    stored as part of the value.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if (-not $map.ContainsKey($r)) {
        continue
    }
    $destIdx0 = $r - $firstRow
    for ($c = 1; $c -le $colCount; $c++) {
        $val = $destValues[$destIdx0, $c - 1]
        $cell = $ws.Cells.Item($r, $firstCol + $c - 1)
        if ($val -eq $null) {
            $cell.Value2 = $null
        } elseif ($val -is [string]) {
            if ($val -eq "") {
                $cell.Value2 = "'"
            } else {
                $cell.Value2 = "'" + $val
            }
        } else {
            $cell.Value2 = $val
        }
    }
}
